# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt -
# Pepino dulce" at row 81. Excel shifts every existing row 81-112 down to
# 82-113 (their data is left untouched), and the new row 81 carries the
# latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(81).Insert()

$ws.Range("A81").Value = 4
$ws.Range("B81").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C81").Value = "Los Lagos"
$ws.Range("D81").Value = 45093
$ws.Range("E81").Value = 10
$ws.Range("F81").Value = 100112043
$ws.Range("G81").Value = "Pepino dulce"
$ws.Range("H81").Value = "Cultivar IV Región"
$ws.Range("I81").Value = "Especial"
$ws.Range("J81").Value = 50
$ws.Range("K81").Value = 20000
$ws.Range("L81").Value = 20000
$ws.Range("M81").Value = 20000
$ws.Range("N81").Value = "`$/bandeja 18 kilos"
$ws.Range("O81").Value = "Provincia de Limarí"
$ws.Range("P81").Value = 1111
$ws.Range("Q81").Value = 18
$ws.Range("R81").Value = "Hortaliza"
